# Populate the previously-empty "Custom Drone 1" sheet with the same header /
# initial-condition rows that already exist on "DJI Phantom 4" (rows 1-3),
# mirroring the data layout used across the other drone sheets in this
# workbook ("Add values to empty sheet for test purpose").

$wb = $excel.ActiveWorkbook

$row1 = @("Time","x","y","z","vx","vy","vz","theta (axe 1)","phi (axe 1)","theta (axe 2)","phi (axe 2)","theta (axe 3)","phi (axe 3)")
$row2 = @("t0","x0","y0","z0","vx0","vy0","vz0","theta0_ax1","phi0_ax1","theta0_ax2","phi0_ax2","theta0_ax3","phi0_ax3")
$row3 = @("t1","x1","y1","z1","vx1","vy1","vz1","theta0_ax2","phi0_ax2","theta0_ax3","phi0_ax3","theta0_ax4","phi0_ax4")

$ws3 = $wb.Worksheets.Item("Custom Drone 1")
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $row1[$i]
    $ws3.Cells.Item(2, $i + 1).Value = $row2[$i]
    $ws3.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Match the column widths that Excel auto-fit for the identical header row on
# the other sheets (best effort - target character widths from the sibling
# "DJI Phantom 4"/"DJI Mavic 3" sheets).
$ws3.Columns.Item(1).ColumnWidth = 4.17
$ws3.Columns.Item(2).ColumnWidth = 2.0
$ws3.Columns.Item(3).ColumnWidth = 2.0
$ws3.Columns.Item(4).ColumnWidth = 2.0
$ws3.Columns.Item(5).ColumnWidth = 3.0
$ws3.Columns.Item(6).ColumnWidth = 3.0
$ws3.Columns.Item(7).ColumnWidth = 2.83
$ws3.Columns.Item(8).ColumnWidth = 10.33
$ws3.Columns.Item(9).ColumnWidth = 8.5
$ws3.Columns.Item(10).ColumnWidth = 10.33
$ws3.Columns.Item(11).ColumnWidth = 8.5
$ws3.Columns.Item(12).ColumnWidth = 10.33
$ws3.Columns.Item(13).ColumnWidth = 8.5

# Update the selection on "DJI Phantom 4" sheet (A1:M2 -> A1:M3) without
# changing which sheet is the active/selected tab.
$ws1 = $wb.Worksheets.Item("DJI Phantom 4")
$ws1.Range("A1:M3").Select()

# Re-activate the "Custom Drone 1" sheet and select its new data range so it
# remains the tab that is shown/selected when the workbook is reopened.
$ws3.Activate()
$ws3.Range("A1:M3").Select()
